$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 217.41667
$ws.Range("I9").Value = 212.33333
$ws.Range("J9").Value = 222.5
$ws.Range("K9").Value = 212.33333
$ws.Range("L9").Value = 222.5
$ws.Range("M9").Value = -43.33332999999999
$ws.Range("N9").Value = -560.5
$ws.Range("H41").Value = 1196.421
$ws.Range("I41").Value = 1125.5
$ws.Range("K41").Value = 1125.5
$ws.Range("M41").Value = -685.5
$ws.Range("H43").Value = 5345.25
$ws.Range("I43").Value = 2224
$ws.Range("K43").Value = 2224
$ws.Range("M43").Value = -2155
$ws.Range("H98").Value = 52640600
$ws.Range("J98").Value = 13080
$ws.Range("L98").Value = 13080
$ws.Range("N98").Value = -16076
$ws.Range("H107").Value = 4191.7144
$ws.Range("I107").Value = 4001.5
$ws.Range("J107").Value = 4267.8
$ws.Range("K107").Value = 4001.5
$ws.Range("L107").Value = 4267.8
$ws.Range("M107").Value = -2081.5
$ws.Range("N107").Value = -8107.8
$ws.Range("H122").Value = 52640600
$ws.Range("J122").Value = 13080
$ws.Range("L122").Value = 39240
$ws.Range("N122").Value = -44140
$ws.Range("H132").Value = 2276.9375
$ws.Range("I132").Value = 1995.4
$ws.Range("K132").Value = 5986.200000000001
$ws.Range("M132").Value = -3456.200000000001
$ws.Range("H138").Value = 2504.875
$ws.Range("I138").Value = 917.5
$ws.Range("J138").Value = 3034
$ws.Range("K138").Value = 2752.5
$ws.Range("L138").Value = 9102
$ws.Range("M138").Value = 2387.5
$ws.Range("N138").Value = -19382
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2311.4
$ws.Range("I45").Value = 2311.4
$ws.Range("K45").Value = 2311.4
$ws.Range("M45").Value = -1934.4
$ws.Range("H54").Value = 24747.5
$ws.Range("J54").Value = 24747.5
$ws.Range("L54").Value = 24747.5
$ws.Range("N54").Value = -26285.5
$ws.Range("H74").Value = 12510189
$ws.Range("I74").Value = 25001644
$ws.Range("K74").Value = 25001644
$ws.Range("M74").Value = -25000770
$ws.Range("H77").Value = 12510189
$ws.Range("I77").Value = 25001644
$ws.Range("K77").Value = 125008220
$ws.Range("M77").Value = -125003852
$ws.Range("H109").Value = 47000
$ws.Range("J109").Value = 47000
$ws.Range("L109").Value = 47000
$ws.Range("N109").Value = -49774
$ws.Range("H122").Value = 1135.091
$ws.Range("I122").Value = 936
$ws.Range("J122").Value = 1666
$ws.Range("K122").Value = 2808
$ws.Range("L122").Value = 4998
$ws.Range("M122").Value = -358
$ws.Range("N122").Value = -9898
$ws.Range("H132").Value = 5391.2983
$ws.Range("I132").Value = 2257.475
$ws.Range("K132").Value = 6772.424999999999
$ws.Range("M132").Value = -4242.424999999999
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2877.5557
$ws.Range("J86").Value = 3633.3333
$ws.Range("L86").Value = 3633.3333
$ws.Range("N86").Value = -5879.3333
$ws.Range("H89").Value = 2877.5557
$ws.Range("J89").Value = 3633.3333
$ws.Range("L89").Value = 18166.6665
$ws.Range("N89").Value = -29398.6665
$ws.Range("H107").Value = 1997.4445
$ws.Range("I107").Value = 1897.7142
$ws.Range("J107").Value = 2346.5
$ws.Range("K107").Value = 1897.7142
$ws.Range("L107").Value = 2346.5
$ws.Range("M107").Value = 22.28580000000011
$ws.Range("N107").Value = -6186.5
$ws.Range("H134").Value = 32344.719
$ws.Range("I134").Value = 1129.5807
$ws.Range("J134").Value = 1000014
$ws.Range("K134").Value = 3388.7421
$ws.Range("L134").Value = 3000042
$ws.Range("M134").Value = -853.7420999999999
$ws.Range("N134").Value = -3005112
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 585448.25
$ws.Range("I31").Value = 1812.4445
$ws.Range("K31").Value = 1812.4445
$ws.Range("M31").Value = -1517.4445
$ws.Range("H34").Value = 585448.25
$ws.Range("I34").Value = 1812.4445
$ws.Range("K34").Value = 1812.4445
$ws.Range("M34").Value = -1610.4445
$ws.Range("H55").Value = 4166.6665
$ws.Range("I55").Value = 4166.6665
$ws.Range("K55").Value = 4166.6665
$ws.Range("M55").Value = -3851.6665
$ws.Range("H58").Value = 2214.3215
$ws.Range("I58").Value = 1999.4348
$ws.Range("J58").Value = 3202.8
$ws.Range("K58").Value = 1999.4348
$ws.Range("L58").Value = 3202.8
$ws.Range("M58").Value = -1796.4348
$ws.Range("N58").Value = -3608.8
$ws.Range("H108").Value = 81477.336
$ws.Range("J108").Value = 81477.336
$ws.Range("L108").Value = 81477.336
$ws.Range("N108").Value = -89157.336
$ws.Range("H132").Value = 2087.923
$ws.Range("I132").Value = 2114.3
$ws.Range("K132").Value = 6342.900000000001
$ws.Range("M132").Value = -3812.900000000001
$ws.Range("H134").Value = 253140.6
$ws.Range("I134").Value = 346510.7
$ws.Range("J134").Value = 6983.091
$ws.Range("K134").Value = 1039532.1
$ws.Range("L134").Value = 20949.273
$ws.Range("M134").Value = -1036997.1
$ws.Range("N134").Value = -26019.273
$ws.Range("H136").Value = 2214.3215
$ws.Range("I136").Value = 1999.4348
$ws.Range("J136").Value = 3202.8
$ws.Range("K136").Value = 5998.3044
$ws.Range("L136").Value = 9608.400000000001
$ws.Range("M136").Value = -3448.3044
$ws.Range("N136").Value = -14708.4
$ws.Range("H139").Value = 97206
$ws.Range("J139").Value = 97206
$ws.Range("L139").Value = 97206
$ws.Range("N139").Value = -107486
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 23605.176
$ws.Range("I2").Value = 132.28572
$ws.Range("K2").Value = 793.71432
$ws.Range("M2").Value = -680.71432
$ws.Range("H132").Value = 2050.0454
$ws.Range("I132").Value = 2144.5
$ws.Range("J132").Value = 1936.7
$ws.Range("K132").Value = 19300.5
$ws.Range("L132").Value = 17430.3
$ws.Range("M132").Value = -16770.5
$ws.Range("N132").Value = -22490.3
$ws.Range("H138").Value = 1806.5714
$ws.Range("I138").Value = 929.2
$ws.Range("K138").Value = 2787.6
$ws.Range("M138").Value = 2352.4
$ws.Range("H140").Value = 302153.8
$ws.Range("I140").Value = 302153.8
$ws.Range("K140").Value = 906461.3999999999
$ws.Range("M140").Value = -901281.3999999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 57500
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 57500
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 57500
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -57812
$ws.Range("H122").Value = 1599.8462
$ws.Range("I122").Value = 1459.8
$ws.Range("J122").Value = 2066.6667
$ws.Range("K122").Value = 4379.4
$ws.Range("L122").Value = 6200.000100000001
$ws.Range("M122").Value = -1929.4
$ws.Range("N122").Value = -11100.0001
$ws.Range("H132").Value = 66669220
$ws.Range("I132").Value = 71430990
$ws.Range("J132").Value = 4395
$ws.Range("K132").Value = 214292970
$ws.Range("L132").Value = 13185
$ws.Range("M132").Value = -214290440
$ws.Range("N132").Value = -18245
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1609.5
$ws.Range("I16").Value = 1328
$ws.Range("K16").Value = 1328
$ws.Range("M16").Value = -1158
$ws.Range("H40").Value = 3607.8235
$ws.Range("I40").Value = 2485.2727
$ws.Range("J40").Value = 5665.8335
$ws.Range("K40").Value = 2485.2727
$ws.Range("L40").Value = 5665.8335
$ws.Range("M40").Value = -2349.2727
$ws.Range("N40").Value = -5937.8335
$ws.Range("H46").Value = 2496.9
$ws.Range("I46").Value = 2395.9
$ws.Range("K46").Value = 2395.9
$ws.Range("M46").Value = -2207.9
$ws.Range("H68").Value = 2515.4666
$ws.Range("I68").Value = 2398.3
$ws.Range("J68").Value = 2749.8
$ws.Range("K68").Value = 2398.3
$ws.Range("L68").Value = 2749.8
$ws.Range("M68").Value = -1649.3
$ws.Range("N68").Value = -4247.8
$ws.Range("H71").Value = 2515.4666
$ws.Range("I71").Value = 2398.3
$ws.Range("J71").Value = 2749.8
$ws.Range("K71").Value = 11991.5
$ws.Range("L71").Value = 13749
$ws.Range("M71").Value = -8247.5
$ws.Range("N71").Value = -21237
$ws.Range("H132").Value = 225068.56
$ws.Range("I132").Value = 202222.4
$ws.Range("J132").Value = 253626.25
$ws.Range("K132").Value = 606667.2
$ws.Range("L132").Value = 760878.75
$ws.Range("M132").Value = -604137.2
$ws.Range("N132").Value = -765938.75
$ws.Range("H136").Value = 48284.152
$ws.Range("I136").Value = 3927.5334
$ws.Range("J136").Value = 108770.45
$ws.Range("K136").Value = 11782.6002
$ws.Range("L136").Value = 326311.35
$ws.Range("M136").Value = -9232.600199999999
$ws.Range("N136").Value = -331411.35
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 80007
$ws.Range("J15").Value = 80007
$ws.Range("L15").Value = 80007
$ws.Range("N15").Value = -80583
$ws.Range("H96").Value = 10999
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("H122").Value = 5473.0415
$ws.Range("I122").Value = 3540.6428
$ws.Range("K122").Value = 10621.9284
$ws.Range("M122").Value = -8171.928400000001
$ws.Range("H132").Value = 3531.6
$ws.Range("I132").Value = 3146.2222
$ws.Range("K132").Value = 9438.6666
$ws.Range("M132").Value = -6908.6666
$ws.Range("H136").Value = 3468.125
$ws.Range("I136").Value = 1948.75
$ws.Range("J136").Value = 4987.5
$ws.Range("K136").Value = 5846.25
$ws.Range("L136").Value = 14962.5
$ws.Range("M136").Value = -3296.25
$ws.Range("N136").Value = -20062.5
